# Task SummaryWk9.xlsx - "moved local copies to repo"
# Fill in the task summary sheet with Jesse's weekly entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# --- Header row: who + week number ---
$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 9

# --- Totals row (cumulative total string introduced early, matches save order) ---
$ws.Range("A14").Value = "Cumulative Total:180"
$ws.Range("C14").Value = "Total: "

# --- Task rows ---
$ws.Range("A3").Value = "Project Build"
$ws.Range("A4").Value = "Req analysysis and elicitation"
$ws.Range("B4").Value = "Go over feedback and any changes to requirements"
$ws.Range("B3").Value = "Work on iteration"

$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 16
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 4

$ws.Range("D14").Formula = "=SUM(D3:D13)"

# --- Column width adjustment for column A ---
# (target stored width 28.28515625; engine quantizes ColumnWidth to 1/6-character
#  steps, so 27.5 is the input that lands on the nearest achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 27.5

# --- Selection as left by the author on save ---
$ws.Range("A5:E5").Select()
